$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-7 and add new rows 8-13 per the
# revised NATMI LR-pair computation (Dr Hou advice).

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vegfc"
$ws.Range("C2").Value = "Vipr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.663313
$ws.Range("H2").Value = 7.989939000000001
$ws.Range("I2").Value = 0.3794306644527501
$ws.Range("J2").Value = 0.3794306644527502
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.02608566666666666
$ws.Range("N2").Value = 0.07825699999999999
$ws.Range("O2").Value = 0.007055522672798636
$ws.Range("P2").Value = 0.007055522672798637
$ws.Range("Q2").Value = 0.06947429514699999
$ws.Range("R2").Value = 0.625268656323
$ws.Range("S2").Value = 0.00267708165580143
$ws.Range("T2").Value = 0.002677081655801431

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vegfc"
$ws.Range("C3").Value = "Vipr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.663313
$ws.Range("H3").Value = 7.989939000000001
$ws.Range("I3").Value = 0.3794306644527501
$ws.Range("J3").Value = 0.3794306644527502
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.6658376666666667
$ws.Range("N3").Value = 1.997513
$ws.Range("O3").Value = 0.1800924934601381
$ws.Range("P3").Value = 0.1800924934601381
$ws.Range("Q3").Value = 1.773334113523
$ws.Range("R3").Value = 15.960007021707
$ws.Range("S3").Value = 0.06833261445653274
$ws.Range("T3").Value = 0.06833261445653276

# Row 4: ECs -> M1
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Vegfc"
$ws.Range("C4").Value = "Vipr2"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.663313
$ws.Range("H4").Value = 7.989939000000001
$ws.Range("I4").Value = 0.3794306644527501
$ws.Range("J4").Value = 0.3794306644527502
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1152503333333333
$ws.Range("N4").Value = 0.345751
$ws.Range("O4").Value = 0.03117234266126738
$ws.Range("P4").Value = 0.03117234266126739
$ws.Range("Q4").Value = 0.306947711021
$ws.Range("R4").Value = 2.762529399189
$ws.Range("S4").Value = 0.01182774268851349
$ws.Range("T4").Value = 0.0118277426885135

# Row 5: ECs -> sCs
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Vegfc"
$ws.Range("C5").Value = "Vipr2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.663313
$ws.Range("H5").Value = 7.989939000000001
$ws.Range("I5").Value = 0.3794306644527501
$ws.Range("J5").Value = 0.3794306644527502
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.890024666666667
$ws.Range("N5").Value = 8.670074
$ws.Range("O5").Value = 0.7816796412057959
$ws.Range("P5").Value = 0.7816796412057959
$ws.Range("Q5").Value = 7.697040265054
$ws.Range("R5").Value = 69.273362385486
$ws.Range("S5").Value = 0.2965932256519025
$ws.Range("T5").Value = 0.2965932256519025

# Row 6: FAPs -> ECs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Vegfc"
$ws.Range("C6").Value = "Vipr2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.178631333333334
$ws.Range("H6").Value = 9.535894000000001
$ws.Range("I6").Value = 0.4528458348143826
$ws.Range("J6").Value = 0.4528458348143827
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.02608566666666666
$ws.Range("N6").Value = 0.07825699999999999
$ws.Range("O6").Value = 0.007055522672798636
$ws.Range("P6").Value = 0.007055522672798637
$ws.Range("Q6").Value = 0.08291671741755555
$ws.Range("R6").Value = 0.7462504567579999
$ws.Range("S6").Value = 0.003195064054815303
$ws.Range("T6").Value = 0.003195064054815303

# Row 7: FAPs -> FAPs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Vegfc"
$ws.Range("C7").Value = "Vipr2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.178631333333334
$ws.Range("H7").Value = 9.535894000000001
$ws.Range("I7").Value = 0.4528458348143826
$ws.Range("J7").Value = 0.4528458348143827
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.6658376666666667
$ws.Range("N7").Value = 1.997513
$ws.Range("O7").Value = 0.1800924934601381
$ws.Range("P7").Value = 0.1800924934601381
$ws.Range("Q7").Value = 2.116452470180223
$ws.Range("R7").Value = 19.048072231622
$ws.Range("S7").Value = 0.08155413554475996
$ws.Range("T7").Value = 0.08155413554475997

# Row 8: FAPs -> M1
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Vegfc"
$ws.Range("C8").Value = "Vipr2"
$ws.Range("D8").Value = "M1"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.178631333333334
$ws.Range("H8").Value = 9.535894000000001
$ws.Range("I8").Value = 0.4528458348143826
$ws.Range("J8").Value = 0.4528458348143827
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1152503333333333
$ws.Range("N8").Value = 0.345751
$ws.Range("O8").Value = 0.03117234266126738
$ws.Range("P8").Value = 0.03117234266126739
$ws.Range("Q8").Value = 0.3663383207104445
$ws.Range("R8").Value = 3.297044886394001
$ws.Range("S8").Value = 0.01411626553556162
$ws.Range("T8").Value = 0.01411626553556163

# Row 9: FAPs -> sCs
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Vegfc"
$ws.Range("C9").Value = "Vipr2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.178631333333334
$ws.Range("H9").Value = 9.535894000000001
$ws.Range("I9").Value = 0.4528458348143826
$ws.Range("J9").Value = 0.4528458348143827
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.890024666666667
$ws.Range("N9").Value = 8.670074
$ws.Range("O9").Value = 0.7816796412057959
$ws.Range("P9").Value = 0.7816796412057959
$ws.Range("Q9").Value = 9.186322959572889
$ws.Range("R9").Value = 82.676906636156
$ws.Range("S9").Value = 0.3539803696792457
$ws.Range("T9").Value = 0.3539803696792458

# Row 10: sCs -> ECs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Vegfc"
$ws.Range("C10").Value = "Vipr2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.177290666666667
$ws.Range("H10").Value = 3.531872
$ws.Range("I10").Value = 0.1677235007328671
$ws.Range("J10").Value = 0.1677235007328671
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.02608566666666666
$ws.Range("N10").Value = 0.07825699999999999
$ws.Range("O10").Value = 0.007055522672798636
$ws.Range("P10").Value = 0.007055522672798637
$ws.Range("Q10").Value = 0.03071041190044444
$ws.Range("R10").Value = 0.276393707104
$ws.Range("S10").Value = 0.001183376962181903
$ws.Range("T10").Value = 0.001183376962181903

# Row 11: sCs -> FAPs
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Vegfc"
$ws.Range("C11").Value = "Vipr2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.177290666666667
$ws.Range("H11").Value = 3.531872
$ws.Range("I11").Value = 0.1677235007328671
$ws.Range("J11").Value = 0.1677235007328671
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.6658376666666667
$ws.Range("N11").Value = 1.997513
$ws.Range("O11").Value = 0.1800924934601381
$ws.Range("P11").Value = 0.1800924934601381
$ws.Range("Q11").Value = 0.7838844704817778
$ws.Range("R11").Value = 7.054960234336001
$ws.Range("S11").Value = 0.03020574345884533
$ws.Range("T11").Value = 0.03020574345884533

# Row 12: sCs -> M1
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Vegfc"
$ws.Range("C12").Value = "Vipr2"
$ws.Range("D12").Value = "M1"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.177290666666667
$ws.Range("H12").Value = 3.531872
$ws.Range("I12").Value = 0.1677235007328671
$ws.Range("J12").Value = 0.1677235007328671
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.1152503333333333
$ws.Range("N12").Value = 0.345751
$ws.Range("O12").Value = 0.03117234266126738
$ws.Range("P12").Value = 0.03117234266126739
$ws.Range("Q12").Value = 0.1356831417635556
$ws.Range("R12").Value = 1.221148275872
$ws.Range("S12").Value = 0.005228334437192265
$ws.Range("T12").Value = 0.005228334437192265

# Row 13: sCs -> sCs
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Vegfc"
$ws.Range("C13").Value = "Vipr2"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.177290666666667
$ws.Range("H13").Value = 3.531872
$ws.Range("I13").Value = 0.1677235007328671
$ws.Range("J13").Value = 0.1677235007328671
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.890024666666667
$ws.Range("N13").Value = 8.670074
$ws.Range("O13").Value = 0.7816796412057959
$ws.Range("P13").Value = 0.7816796412057959
$ws.Range("Q13").Value = 3.402399066503111
$ws.Range("R13").Value = 30.621591598528
$ws.Range("S13").Value = 0.1311060458746476
$ws.Range("T13").Value = 0.1311060458746476
